# Update the "Last Updated" timestamp on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 10:41 AM"

# Insert a new top row of data into the "Stock List" sheet: a new entry
# (MIDWESTLTD) is prepended at row 2, shifting every existing data row down
# by one (the former last row falls off the bottom of the table).
$ws = $wb.Worksheets.Item("Stock List")

for ($r = 75; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    $ws.Range("B$dst").Value = $ws.Range("B$src").Value2
    $ws.Range("C$dst").Value = $ws.Range("C$src").Value2
    $ws.Range("D$dst").Value = $ws.Range("D$src").Value2
    $ws.Range("E$dst").Value = $ws.Range("E$src").Value2
    $ws.Range("H$dst").Value = $ws.Range("H$src").Value2
}

$ws.Range("B2").Value = "MIDWESTLTD"
$ws.Range("C2").Value = "MIDWESTLTD"
$ws.Range("D2").Value = 1117.2
$ws.Range("E2").Value = -1.4032
$ws.Range("H2").Value = 4039.8864
